$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.412.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.849.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'240.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'0.6299"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07687"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").Value = "'0.2943"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "'24.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("D11").Value = "'0.07748"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "'1.847.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "'5.027"
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = "  +7.98%  "
$ws.Range("D15").Value = "'0.6810"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "'83.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "'2.096.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "'6.152"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'29.422.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'229.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'7.456"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'157.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").Value = "'8.387"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "'17.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'1.316"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.08%  "
$ws.Range("D30").Value = "'1.467"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "'0.05718"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "'1.851"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "'0.7076"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'2.780"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Value = "'1.224.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("D41").Value = "'6.445"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.48%  "
$ws.Range("D42").Value = "'0.9114"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'101.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'66.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000121"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.98%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.143"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").Value = "'0.4026"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.033"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.691"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1125"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
